# ToDo.xlsx edit: rename sheets, move active tab, add new checklist rows.
$wb = $excel.ActiveWorkbook

# --- Sheet renames -----------------------------------------------------
# Tab order in the workbook is [Sheet2, Sheet1] (r:id rId1 -> physical
# sheet1.xml, rId2 -> physical sheet2.xml). Rename them to the new titles.
$wsStatus = $wb.Worksheets.Item(1)   # was "Sheet2" -> "Status"
$wsFix    = $wb.Worksheets.Item(2)   # was "Sheet1" -> "Have to fix"
$wsStatus.Name = "Status"
$wsFix.Name = "Have to fix"

# --- "Status" sheet: re-split the last checklist entry ------------------
# Row 43 gains a "not completed" cell in column B (same text/style that
# used to live alone in row 44); row 44 then becomes the old row 45
# ("none exist"), and the sheet shrinks by one row.
$wsStatus.Range("A44").Copy($wsStatus.Range("B43"))
$wsStatus.Range("A45").Copy($wsStatus.Range("A44"))
$wsStatus.Rows.Item(45).Delete()

# --- "Have to fix" sheet: add a new to-do item --------------------------
$wsFix.Range("A10").Value = "Correct meaning key value of all components"

# --- View state ----------------------------------------------------------
# Selection on "Have to fix" moves to A11 (it is no longer the active tab,
# so set this first - selecting on "Status" below switches the active tab).
$wsFix.Activate() | Out-Null
$wsFix.Range("A11").Select() | Out-Null

# Move the scroll/selection on "Status" further down, select cell B44, and
# leave "Status" as the active tab (it is now tab index 0).
$wsStatus.Activate() | Out-Null
$wsStatus.Range("B44").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
